# Update "want to go" counts (column F) on both the "展览" (sheet 1) and
# "全部类型" (sheet 4) worksheets, matching the refreshed scrape output.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAll        = $wb.Worksheets.Item("全部类型")

# Row -> new value updates for the "展览" sheet
$exhibitionUpdates = @{
    2  = 14950
    3  = 18780
    13 = 56
    15 = 213
    22 = 7810
    28 = 5995
    33 = 274
    34 = 5365
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Row -> new value updates for the "全部类型" sheet
$allTypesUpdates = @{
    2  = 14950
    3  = 18780
    13 = 56
    15 = 213
    23 = 7810
    31 = 5995
    36 = 274
    37 = 5365
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAll.Range("F$row").Value = $allTypesUpdates[$row]
}
